# Generate Report for Handback
# Updates timestamps / priority produced by a (re)generated handback report.

$wb = $excel.ActiveWorkbook

# --- Sheet "Overview" ---
$wsOverview = $wb.Worksheets.Item("Overview")
# "Latest HO Xliff Generate Date" for the 2fed587e... file (rows 2 and 3 share the value)
$wsOverview.Range("G2").Value = "2016-09-06 14:25:23"
$wsOverview.Range("G3").Value = "2016-09-06 14:25:23"

# --- Sheet "zh-cn" ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
# Priority changed from "ht" (human translation) to "mt" (machine translation)
$wsZhCn.Range("E2").Value = "mt"
$wsZhCn.Range("E3").Value = "mt"
# Correspond Handoff Datetime
$wsZhCn.Range("H2").Value = "2016-09-06 14:24:59"
$wsZhCn.Range("H3").Value = "2016-09-06 14:24:59"
# Correspond Handback DateTime
$wsZhCn.Range("K2").Value = "2016-09-06 14:25:48"
$wsZhCn.Range("K3").Value = "2016-09-06 14:25:48"

# --- Sheet "de-de" ---
$wsDeDe = $wb.Worksheets.Item("de-de")
# Priority changed from "ht" (human translation) to "mt" (machine translation)
$wsDeDe.Range("E2").Value = "mt"
$wsDeDe.Range("E3").Value = "mt"
# Correspond Handoff Datetime (same underlying value as Overview's generate date)
$wsDeDe.Range("H2").Value = "2016-09-06 14:25:23"
$wsDeDe.Range("H3").Value = "2016-09-06 14:25:23"
# Correspond Handback DateTime
$wsDeDe.Range("K2").Value = "2016-09-06 14:25:56"
$wsDeDe.Range("K3").Value = "2016-09-06 14:25:56"
